# Applies the "ADDITIONAL SCRAPING" edit:
#  - Splits/augments the single "ODI Batting" sheet into 3 sheets:
#      1) "Player Info"        (new)
#      2) "ODI Batting"        (existing, column D repurposed to MATCH_CODE)
#      3) "ODI Batting Extra"  (new)
#
# NOTE: worksheet COM handles in this runtime are positional, so we
# re-fetch each worksheet by name right before using it instead of
# holding on to a reference across structural operations (Add/rename).

$wb = $excel.ActiveWorkbook

# --- locate the existing sheet, then insert the two new sheets -----------------
$battingOrig = $wb.Worksheets.Item(1)
$battingOrigName = $battingOrig.Name

$extra = $wb.Worksheets.Add($null, $wb.Worksheets.Item($battingOrigName))
$extra.Name = "ODI Batting Extra"

$info = $wb.Worksheets.Add($wb.Worksheets.Item($battingOrigName), $null)
$info.Name = "Player Info"

# =================================================================================
# 1) "Player Info" sheet
# =================================================================================
$wsInfo = $wb.Worksheets.Item("Player Info")

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($col = 1; $col -le $infoHeaders.Length; $col++) {
    $cell = $wsInfo.Cells.Item(1, $col)
    $cell.Value = $infoHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$wsInfo.Cells.Item(2, 1).NumberFormat = "@"
$wsInfo.Cells.Item(2, 1).Value = "4769"
$wsInfo.Cells.Item(2, 2).Value = "Shubman Gill"
$wsInfo.Cells.Item(2, 3).Value = "Right Handed"
$wsInfo.Cells.Item(2, 4).Value = "Right Arm Off Break"

$wsInfo.Range("A1").Select()

# =================================================================================
# 2) "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# =================================================================================
$wsBatting = $wb.Worksheets.Item("ODI Batting")

$wsBatting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @(
    "4248", "4249", "4437", "4621", "4623", "4624", "4637", "4640", "4643", "4656",
    "4657", "4658", "4669", "4673", "4676", "4687", "4689", "4691", "4692", "4695",
    "4697", "4725", "4728", "4732"
)

for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $wsBatting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
}

# =================================================================================
# 3) "ODI Batting Extra" sheet
# =================================================================================
$wsExtra = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($col = 1; $col -le $extraHeaders.Length; $col++) {
    $cell = $wsExtra.Cells.Item(1, $col)
    $cell.Value = $extraHeaders[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("4623", $null, $null, $null, $null, "NO"),
    @("4624", 2, "7", "2", "43.56%", "YES"),
    @("4637", 2, "10", "1", "42.71%", "NO"),
    @("4640", 3, "6", "0", "19.76%", "NO"),
    @("4643", $null, $null, $null, $null, "NO"),
    @("4656", $null, $null, $null, $null, "NO"),
    @("4657", 2, "4", "0", "9.93%", "NO"),
    @("4658", 2, "8", "0", "46.67%", "NO"),
    @("4669", 2, "1", "3", "16.34%", "NO"),
    @("4673", $null, $null, $null, $null, "NO"),
    @("4676", 2, "2", "0", "5.94%", "NO"),
    @("4687", $null, $null, $null, $null, "NO"),
    @("4689", 2, "5", "0", "9.59%", "NO"),
    @("4691", $null, $null, $null, $null, "NO"),
    @("4692", $null, $null, $null, $null, "NO"),
    @("4695", 2, "6", "0", "36.04%", "NO"),
    @("4697", 2, "13", "5", "29.09%", "NO"),
    @("4725", 2, "3", "0", "10.47%", "NO"),
    @("4728", 2, "0", "0", $null, "NO"),
    @("4732", 2, "4", "1", "14.92%", "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $row = $i + 2
    $values = $extraRows[$i]

    $codeCell = $wsExtra.Cells.Item($row, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $values[0]

    # BATTING_POSITION (numeric, left untouched/blank when not applicable)
    if ($null -ne $values[1]) {
        $wsExtra.Cells.Item($row, 2).Value = $values[1]
    }

    # NUM_4 / NUM_6 (text-typed, left untouched/blank when not applicable)
    for ($col = 3; $col -le 4; $col++) {
        $v = $values[$col - 1]
        if ($null -ne $v) {
            $cell = $wsExtra.Cells.Item($row, $col)
            $cell.NumberFormat = "@"
            $cell.Value = $v
        }
    }

    # PERCENT_RUNS_OF_TOTAL (text-typed, left untouched/blank when not applicable)
    if ($null -ne $values[4]) {
        $pctCell = $wsExtra.Cells.Item($row, 5)
        $pctCell.NumberFormat = "@"
        $pctCell.Value = $values[4]
    }

    $wsExtra.Cells.Item($row, 6).Value = $values[5]
}

$wsExtra.Range("A1").Select()

# --- make sure the original sheet stays the active one on open -----------------
$wb.Worksheets.Item("ODI Batting").Activate()
